$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The "Baz chan" / bookmark("_GoBack") / "ges" paragraph needs to become two
# paragraphs:
#   1) "Baz changes"                                        (no bookmark)
#   2) "Andrew Mckay ID 22142955 changes have been made."    (bookmark at end)
# plus one extra blank paragraph appended right after that second paragraph.
# ---------------------------------------------------------------------------

# Locate the paragraph that currently reads "Baz chan" + "ges" (split across
# two runs around the _GoBack bookmark) instead of hard-coding its index.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Baz chan*") {
        $targetIndex = $i
        break
    }
}

# The existing _GoBack bookmark sits inside that paragraph; it needs to move
# to the new paragraph we are about to create, so drop it for now.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# Rewrite the paragraph's text (excluding the trailing paragraph mark) so the
# two runs ("Baz chan" / "ges") collapse into a single "Baz changes" run.
# Doing this in two steps (temp text, then final text) forces the engine to
# actually rebuild the run instead of treating an unchanged-text assignment
# as a no-op that leaves the original run split in place.
$p = $d.Paragraphs.Item($targetIndex)
$r = $p.Range
$target = $d.Range($r.Start, $r.End - 1)
$target.Text = "TEMP_PLACEHOLDER"
$p = $d.Paragraphs.Item($targetIndex)
$r = $p.Range
$target = $d.Range($r.Start, $r.End - 1)
$target.Text = "Baz changes"

# Insert a brand new paragraph right after it.
$p = $d.Paragraphs.Item($targetIndex)
$p.Range.InsertParagraphAfter()
$newIndex = $targetIndex + 1

# Fill the new paragraph with its text (again excluding the paragraph mark),
# but append a one-character filler first. Adding the bookmark exactly at
# "paragraph end - 1" is unreliable in this host when that position lines up
# with a run/paragraph boundary, so we park the bookmark one character
# earlier (a safe, non-boundary offset) and then delete the filler character
# afterwards - the bookmark (anchored before the deleted character) ends up
# sitting correctly at the end of the real text, right before the paragraph
# mark.
$newPara = $d.Paragraphs.Item($newIndex)
$r2 = $newPara.Range
$body = $d.Range($r2.Start, $r2.End - 1)
$body.Text = "Andrew Mckay ID 22142955 changes have been made.X"

$newPara = $d.Paragraphs.Item($newIndex)
$r3 = $newPara.Range
$bookmarkPos = $r3.End - 2
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$newPara = $d.Paragraphs.Item($newIndex)
$r4 = $newPara.Range
$fillerRange = $d.Range($r4.End - 2, $r4.End - 1)
$fillerRange.Text = ""

# Add one more blank paragraph right after the paragraph we just created.
$newPara = $d.Paragraphs.Item($newIndex)
$newPara.Range.InsertParagraphAfter()
